$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativação date (row 8, cols B & C) ---
# Plain assignment of "01/01/2022" gets auto-parsed by Excel as a date
# serial (and pulls in a new number-format/style), so instead stage the
# text in a scratch cell via a formula (so it's typed as text, not a
# date), copy it, and paste-special *values only* into the target cells.
# That keeps the original style (General number format) on B8/C8 intact.
$ws.Range("Z1").Formula = "=""01/01/2022"""
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").ClearContents()

# --- Programa resumido (row 16, cols B & C) ---
$novoPrograma = "Estática de Partículas. Estática de Corpos Rígidos. Equilíbrio de Corpos Rígidos. Análise de Estruturas."
$ws.Range("B16").Value = $novoPrograma
$ws.Range("C16").Value = $novoPrograma

# --- Método (row 21, cols B & C) ---
$novoMetodo = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"
$ws.Range("B21").Value = $novoMetodo
$ws.Range("C21").Value = $novoMetodo

# --- Critério (row 22, cols B & C) ---
$novoCriterio = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("B22").Value = $novoCriterio
$ws.Range("C22").Value = $novoCriterio

# --- Norma de recuperação (row 23, cols B & C) ---
$novaNorma = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("B23").Value = $novaNorma
$ws.Range("C23").Value = $novaNorma

# --- Bibliografia (row 24, cols B & C) ---
$novaBiblio = "1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p."
$ws.Range("B24").Value = $novaBiblio
$ws.Range("C24").Value = $novaBiblio
